# Apply the 2021 divisional-round (simulated season) stat updates to the
# Rushing and Receiving sheets of the Rams "Players Data" workbook.

$wb = $excel.ActiveWorkbook

$wsRushing   = $wb.Worksheets.Item("Rushing")
$wsReceiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet -----------------------------------------------------
# Columns: B=Name, C=1DATT, D=2DATT, E=3DATT, F=RZATT

# Row 2 - M.Stafford
$wsRushing.Range("D2").Value = 8
$wsRushing.Range("E2").Value = 10
$wsRushing.Range("F2").Value = 7

# Row 4 - S.Michel
$wsRushing.Range("E4").Value = 34

# Row 7 - C.Akers
$wsRushing.Range("C7").Value = 31
$wsRushing.Range("D7").Value = 14
$wsRushing.Range("F7").Value = 11

# Row 9 - V.Jefferson
$wsRushing.Range("D9").Value = 1

# --- Receiving sheet -----------------------------------------------------
# Columns: B=Name, C=Short Target, D=Short Comp, E=Deep Target,
#          F=Deep Comp, G=RZ Target, H=RZ Comp

# Row 2 - S.Michel
$wsReceiving.Range("C2").Value = 49
$wsReceiving.Range("D2").Value = 33

# Row 3 - C.Akers
$wsReceiving.Range("C3").Value = 6
$wsReceiving.Range("D3").Value = 6

# Row 4 - C.Kupp
$wsReceiving.Range("C4").Value = 146
$wsReceiving.Range("D4").Value = 114
$wsReceiving.Range("E4").Value = 51
$wsReceiving.Range("F4").Value = 32
$wsReceiving.Range("G4").Value = 36
$wsReceiving.Range("H4").Value = 25

# Row 5 - V.Jefferson
$wsReceiving.Range("C5").Value = 64
$wsReceiving.Range("D5").Value = 39
$wsReceiving.Range("G5").Value = 16
$wsReceiving.Range("H5").Value = 5

# Row 7 - O.Beckham
$wsReceiving.Range("C7").Value = 78
$wsReceiving.Range("D7").Value = 64
$wsReceiving.Range("E7").Value = 31
$wsReceiving.Range("F7").Value = 13
$wsReceiving.Range("G7").Value = 24
$wsReceiving.Range("H7").Value = 13

# Row 8 - K.Blanton
$wsReceiving.Range("C8").Value = 5
$wsReceiving.Range("D8").Value = 3
$wsReceiving.Range("G8").Value = 2
$wsReceiving.Range("H8").Value = 1

# Row 9 - T.Higbee
$wsReceiving.Range("C9").Value = 85
$wsReceiving.Range("D9").Value = 63
$wsReceiving.Range("E9").Value = 11
$wsReceiving.Range("F9").Value = 5
